$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (row 1) with two new columns P and Q,
# matching the formatting of the existing header cell O1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# For each data row (2-25): swap values in columns I/K and M/O,
# and add new columns P and Q with value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new, value 2
    $ws.Cells.Item($r, 17).Value = 2   # Q: new, value 2
}
